$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01702384718974029
$ws.Range("C2").Value = 0.2168581401471243

$ws.Range("B3").Value = 0.06596560242209266
$ws.Range("C3").Value = 0.2517934647325251

$ws.Range("B4").Value = 0.8696459864049334
$ws.Range("C4").Value = 0.4923520739298654

$ws.Range("B5").Value = 0.993363207328594
$ws.Range("C5").Value = 0.4579360745264823

$ws.Range("B6").Value = 0.9864096432805773
$ws.Range("C6").Value = 0.7648136320469063

$ws.Range("B7").Value = 0.9495099268375535
$ws.Range("C7").Value = 0.3372410961695695

$ws.Range("B8").Value = 0.00897695541381836
$ws.Range("C8").Value = 0.1923011207580566

$ws.Range("B9").Value = 0.1949137773170574
$ws.Range("C9").Value = 0.2513584432831898

$ws.Range("B10").Value = 0.6069232130487306
$ws.Range("C10").Value = 0.4097708737640559
